$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 28, shifting rows 28:121 down to 29:122
$ws.Rows.Item(28).Insert()

# Fill in the constant / repeated columns for the new row 28 (same as other data rows)
$ws.Cells.Item(28, 1).Value = 8
$ws.Cells.Item(28, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(28, 3).Value = "Coquimbo"
$ws.Cells.Item(28, 4).Value = 44607
$ws.Cells.Item(28, 5).Value = 4
$ws.Cells.Item(28, 6).Value = 100112040
$ws.Cells.Item(28, 7).Value = "Cilantro"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 2400
$ws.Cells.Item(28, 11).Value = 2300
$ws.Cells.Item(28, 12).Value = 2500
$ws.Cells.Item(28, 13).Value = 2400
$ws.Cells.Item(28, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(28, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(28, 16).Value = 1600
$ws.Cells.Item(28, 17).Value = 1.5
$ws.Cells.Item(28, 18).Value = "Hortaliza"
